$d = $word.ActiveDocument

# 1. Merge "Dundee University " + "TouAR" into a single run "Dundee University TouAR"
#    (removes the spell-check-induced run split; the proofErr markers go away too
#    since we are replacing the whole visible text span).
$d.Content.Find.Execute("Dundee University TouAR", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Dundee University TouAR", 2) | Out-Null

# 2. Replace the placeholder "insert approval number from decision letter" with the
#    real approval number, and strip the bold/italic/highlight formatting it had.
$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Text = "insert approval number from decision letter"
$range.Find.Replacement.ClearFormatting()
$range.Find.Replacement.Text = "UOD-SSEREC-DoC-UG-2019-008"
$range.Find.Execute("insert approval number from decision letter", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "UOD-SSEREC-DoC-UG-2019-008", 2) | Out-Null

# Now strip the bold/italic/highlight formatting from the replaced run so it
# reads as plain text.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Text = "UOD-SSEREC-DoC-UG-2019-008"
if ($findRange.Find.Execute()) {
    $findRange.Font.Bold = 0
    $findRange.Font.Italic = 0
    $findRange.HighlightColorIndex = 0
}

# 3. Remove the "_GoBack" bookmark near "Do I have to take part?"
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}
